$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dash = [char]0x2014

# --- Fill in "Box 1" package details for Justice Oladeji (row 2) ---
$ws.Range("C2").Value = "Box 1"
$ws.Range("D2").Value = "INCH"
$ws.Range("E2").Value = 23
$ws.Range("F2").Value = 33
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = 7.7
$ws.Range("I2").Formula = "=E2*2.54*F2*2.54*G2*2.54/6000"
$ws.Range("J2").Value = "Yes"
$ws.Range("K2").Formula = '=IF(J2="No", 13*MAX(H2,I2), IF(J2="Yes", 14*MAX(H2,I2), "Invalid Input' + $dash + 'Yes or No"))'

# --- Fill in "Box 2" package details for Uzoma Emah (row 9) ---
$ws.Range("C9").Value = "Box 2"
$ws.Range("D9").Value = "INCH"
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = 2
$ws.Range("I9").Formula = "=E9*2.54*F9*2.54*G9*2.54/6000"
$ws.Range("J9").Value = "No"
$ws.Range("K9").Formula = '=IF(J9="No", 13*MAX(H9,I9), IF(J9="Yes", 14*MAX(H9,I9), "Invalid Input' + $dash + 'Yes or No"))'

# --- Fill in "Box 3" package details, sharing row 10 (Uzoma's Email row) ---
$ws.Range("C10").Value = "Box 3"
$ws.Range("D10").Value = "INCH"
$ws.Range("E10").Value = 30
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = 3
$ws.Range("I10").Formula = "=E10*2.54*F10*2.54*G10*2.54/6000"
$ws.Range("J10").Value = "Yes"
$ws.Range("K10").Formula = '=IF(J10="No", 13*MAX(H10,I10), IF(J10="Yes", 14*MAX(H10,I10), "Invalid Input' + $dash + 'Yes or No"))'

# --- Apply the bold "package label" formatting (matching the old Box-1 label style)
#     to the newly-filled package-label cells, pulling the format from the
#     existing "Box 1" label cell (C23) before it is removed below.
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Update "Total Cost" formulas to sum the package costs now present ---
$ws.Range("B7").Formula = "=SUM(K2:K2)"
$ws.Range("B14").Formula = "=SUM(K9:K10)"

# --- Remove the now-redundant "John Doe" and "Collins Emenike" blocks (rows 16-29) ---
$ws.Range("A16:K29").EntireRow.Delete() | Out-Null
